$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the test case numbers for "Spint(39) - Day 5" section (rows 27-29):
# created test cases for api_signup_password_field and executed
# api_signup_username and symlex_portocol_settings
$ws.Range("C27").Value = 912
$ws.Range("C28").Value = 1115
$ws.Range("C29").Value = 636

# Update the view: scroll so row 15 is the top-left row, and select C29
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C29").Select()
